$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the built-in "Normalny" cell style to "Normal" (Polish -> English
#        locale normalization seen when the workbook is re-saved from EN Excel). ---
foreach ($st in $wb.Styles) {
    if ($st.Name -eq "Normalny") {
        $st.Name = "Normal"
    }
}

# --- 2. Append the new data rows (78-106) coming from the refreshed
#        sullygnome / supremacy / packages export. ---

# Seed the formatting (date style on col A, time style on col B) for the new
# rows by copying down the formats already used on the last existing row,
# exactly like a user dragging the fill handle down before typing new data.
$ws.Range("A77:D77").Copy()
$ws.Range("A78:D106").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$newRows = @(
    @(78, 45094, 0.72203703703703714, 81200, 1620),
    @(79, 45096, 0.42538194444444444, 81201, 1620),
    @(80, 45097, 0.71953703703703698, 81203, 1620),
    @(81, 45098, 0.89884259259259258, 81219, 1610),
    @(82, 45099, 0.67739583333333331, 81221, 1610),
    @(83, 45101, 0.52664351851851854, 81218, 1610),
    @(84, 45102, 0.46377314814814818, 81256, 1610),
    @(85, 45104, 0.47160879629629626, 81288, 1620),
    @(86, 45105, 0.47478009259259263, 81310, 1620),
    @(87, 45106, 0.76927083333333324, 81300, 1620),
    @(88, 45107, 0.64849537037037031, 81371, 1620),
    @(89, 45108, 0.68618055555555557, 81495, 1620),
    @(90, 45109, 0.60706018518518523, 81495, 1620),
    @(91, 45110, 0.47806712962962966, 81494, 1620),
    @(92, 45111, 0.38723379629629634, 81494, 1620),
    @(93, 45112, 0.45993055555555556, 81498, 1620),
    @(94, 45113, 0.44883101851851853, 81507, 1620),
    @(95, 45114, 0.48599537037037038, 81618, 1620),
    @(96, 45115, 0.46984953703703702, 81618, 1620),
    @(97, 45116, 0.46643518518518517, 81620, 1620),
    @(98, 45117, 0.46270833333333333, 81620, 1620),
    @(99, 45118, 0.5430787037037037, 81620, 1620),
    @(100, 45119, 0.46430555555555553, 81610, 1620),
    @(101, 45120, 0.50746527777777783, 81612, 1620),
    @(102, 45121, 0.58672453703703698, 81613, 1630),
    @(103, 45122, 0.51321759259259259, 81615, 1630),
    @(104, 45123, 0.5237384259259259, 81615, 1630),
    @(105, 45124, 0.50491898148148151, 81615, 1630),
    @(106, 45125, 0.47068287037037032, 81637, 1630)
)

foreach ($r in $newRows) {
    $rowIndex = [int]$r[0]
    $ws.Cells.Item($rowIndex, 1).Value = $r[1]
    $ws.Cells.Item($rowIndex, 2).Value = $r[2]
    $ws.Cells.Item($rowIndex, 3).Value = $r[3]
    $ws.Cells.Item($rowIndex, 4).Value = $r[4]
}

# --- 3. Move the view/selection down to the new last row, as it would be
#        after scrolling through and entering the freshly appended data
#        (scrolls the window so row 86 is the top visible row, then selects
#        the empty cell right after the new data, matching how Excel leaves
#        the cursor after typing in the last row of a pasted/typed block). ---
$excel.Goto($ws.Range("A86"), $true)
$ws.Range("A107").Select()

# --- 4. Best-effort: nudge the application window back towards the
#        top-left of the screen (matches the smaller xWindow/yWindow seen
#        in the diff's workbookView once re-saved on the new machine). ---
$win = $excel.ActiveWindow
$win.Left = 120
$win.Top = 60
